$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16; existing rows 16-25 shift down to 17-26.
$ws.Rows.Item(16).Insert()

# Populate the new row 16 with the new weekly entry.
$ws.Range("A16").Value = 10
$ws.Range("B16").Value = "Vega Modelo de Temuco"
$ws.Range("C16").Value = "La Araucanía"
$ws.Range("D16").Value = 44845
$ws.Range("E16").Value = 9
$ws.Range("F16").Value = 100112036
$ws.Range("G16").Value = "Caigua"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 20
$ws.Range("K16").Value = 16000
$ws.Range("L16").Value = 16000
$ws.Range("M16").Value = 16000
$ws.Range("N16").Value = "$/caja 15 kilos"
$ws.Range("O16").Value = "Región de Arica y Parinacota"
$ws.Range("P16").Value = 1067
$ws.Range("Q16").Value = 15
$ws.Range("R16").Value = "Hortaliza"
